# DP-2544 - add different dataset to second tab and rename tabs
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename tabs ---
$ws1.Name = "Agents Information"
$ws2.Name = "Countries"

# --- Sheet1 ("Agents Information"): widen a couple columns, move selection ---
$ws1.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws1.Columns.Item(4).ColumnWidth = 17.666666666666668

$ws1.Select()
$ws1.Range("K15").Select()

# --- Sheet2 ("Countries"): replace dataset entirely ---
$ws2.Select()
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "id"
$ws2.Range("B1").Value = "country"
$ws2.Range("C1").Value = "Addres"

$rows = @(
    @(1, "USA",    "101 Independecne Avenue"),
    @(1, "USA",    "101 Independecne Avenue"),
    @(2, "Canada", "2796 Reserve St"),
    @(2, "Canada", "2796 Reserve St"),
    @(3, "USA",    "15205 North Kierland"),
    @(3, "USA",    "15205 North Kierland"),
    @(4, "Poland", "ul. Gdańca Pawła 36"),
    @(4, "Poland", "ul. Gdańca Pawła 36")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws2.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws2.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws2.Range("A2:A9").NumberFormat = "0"

$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 28.166666666666668
$ws2.Columns.Item(5).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 17

$ws2.Range("G13:G14").Select()
$ws2.Range("G14").Activate()

$ws2.PageSetup.Orientation = 1

Write-Host "done"
